$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Back to Contents" banner row (old row 1) was removed and every row
# below it shifted up by one. A native row delete keeps cell styles/merges
# consistent for everything that isn't a hyperlink anchor.
$ws.Rows("1").Delete()

# This runtime does not re-anchor the worksheet's Hyperlinks collection when
# rows are deleted, so rebuild the three remaining hyperlinks (the old
# "Back to Contents" one pointed at B1 and is simply gone) at their new,
# shifted addresses: B35->B34, B38->B37, B37->B36.
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("B34"), "https://oe.cd/AFDD-2023") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B37"), "https://github.com/AfDDAnnex/AfDDDDAf2023/raw/main/AfDD_2023_Stats%20(CSV).zip") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B36"), "https://github.com/AfDDAnnex/AfDDDDAf2023/raw/main/AfDD_2023_Stats_by_year.xlsx") | Out-Null

# Hyperlinks.Add() stamps the generic "Hyperlink" look on its anchor cell,
# which wipes the bespoke bold/italic/size-12 formatting B34 ("List of
# sources" link text) carried before the shift. Put it back; B36/B37 keep
# the plain Hyperlink style Add() already gave them (matches their original
# formatting).
$b34 = $ws.Range("B34")
$b34.Style = "Hyperlink"
$b34.Font.Bold = $true
$b34.Font.Italic = $true
$b34.Font.Size = 12

# Match the recorded selection state: the whole first row is selected after
# the header row disappears.
$ws.Rows("1").Select() | Out-Null
